$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 325.75
$ws.Range("J19").Value = 325.75
$ws.Range("L19").Value = 325.75
$ws.Range("N19").Value = -675.75

$ws.Range("H32").Value = 13892055
$ws.Range("J32").Value = 4332.3335
$ws.Range("L32").Value = 4332.3335
$ws.Range("N32").Value = -4984.3335

$ws.Range("H112").Value = 1568.1111
$ws.Range("J112").Value = 1609.7693
$ws.Range("L112").Value = 4829.3079
$ws.Range("N112").Value = -7045.3079

$ws.Range("H138").Value = 1723.579
$ws.Range("I138").Value = 983.2
$ws.Range("J138").Value = 4500
$ws.Range("K138").Value = 2949.6
$ws.Range("L138").Value = 13500
$ws.Range("M138").Value = 2190.4
$ws.Range("N138").Value = -23780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 147992.55
$ws.Range("I74").Value = 223484.4
$ws.Range("K74").Value = 223484.4
$ws.Range("M74").Value = -222610.4

$ws.Range("H77").Value = 147992.55
$ws.Range("I77").Value = 223484.4
$ws.Range("K77").Value = 1117422
$ws.Range("M77").Value = -1113054

$ws.Range("H88").Value = 4042.2
$ws.Range("J88").Value = 7399.75
$ws.Range("L88").Value = 7399.75
$ws.Range("N88").Value = -8211.75

$ws.Range("H91").Value = 4042.2
$ws.Range("J91").Value = 7399.75
$ws.Range("L91").Value = 7399.75
$ws.Range("N91").Value = -10207.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 5318.5
$ws.Range("I33").Value = 290
$ws.Range("J33").Value = 6994.6665
$ws.Range("K33").Value = 1740
$ws.Range("L33").Value = 41967.999
$ws.Range("M33").Value = -1457
$ws.Range("N33").Value = -42533.999

$ws.Range("H40").Value = 275.33334
$ws.Range("J40").Value = 461
$ws.Range("L40").Value = 1844
$ws.Range("N40").Value = -1982

$ws.Range("H44").Value = 4541
$ws.Range("J44").Value = 5125.4287
$ws.Range("L44").Value = 15376.2861
$ws.Range("N44").Value = -16172.2861

$ws.Range("H54").Value = 4500
$ws.Range("I54").Value = 2500
$ws.Range("K54").Value = 7500
$ws.Range("M54").Value = -6941

$ws.Range("H68").Value = 4766891.5
$ws.Range("J68").Value = 10006898
$ws.Range("L68").Value = 30020694
$ws.Range("N68").Value = -30022316

$ws.Range("H69").Value = 1624
$ws.Range("I69").Value = 1624
$ws.Range("K69").Value = 4872
$ws.Range("M69").Value = -4061

$ws.Range("H71").Value = 4766891.5
$ws.Range("J71").Value = 10006898
$ws.Range("L71").Value = 90062082
$ws.Range("N71").Value = -90070194

$ws.Range("H72").Value = 1624
$ws.Range("I72").Value = 1624
$ws.Range("K72").Value = 14616
$ws.Range("M72").Value = -10560

$ws.Range("H80").Value = 7250
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 7250
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 21750
$ws.Range("N80").Value = -23622
$ws.Range("M80").ClearContents()

$ws.Range("H82").Value = 11249.917
$ws.Range("J82").Value = 11818.182
$ws.Range("L82").Value = 35454.546
$ws.Range("N82").Value = -36266.546

$ws.Range("H83").Value = 7250
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 7250
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 65250
$ws.Range("N83").Value = -74610
$ws.Range("M83").ClearContents()

$ws.Range("H85").Value = 11249.917
$ws.Range("J85").Value = 11818.182
$ws.Range("L85").Value = 35454.546
$ws.Range("N85").Value = -38262.546

$ws.Range("H86").Value = 2736.25
$ws.Range("I86").Value = 470
$ws.Range("J86").Value = 3491.6667
$ws.Range("K86").Value = 1410
$ws.Range("L86").Value = 10475.0001
$ws.Range("M86").Value = -224
$ws.Range("N86").Value = -12847.0001

$ws.Range("H87").Value = 1500
$ws.Range("I87").Value = 1500
$ws.Range("K87").Value = 4500
$ws.Range("M87").Value = -3252

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H89").Value = 2736.25
$ws.Range("I89").Value = 470
$ws.Range("J89").Value = 3491.6667
$ws.Range("K89").Value = 4230
$ws.Range("L89").Value = 31425.0003
$ws.Range("M89").Value = 1698
$ws.Range("N89").Value = -43281.0003

$ws.Range("H90").Value = 1500
$ws.Range("I90").Value = 1500
$ws.Range("K90").Value = 13500
$ws.Range("M90").Value = -7260

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H127").Value = 4500
$ws.Range("J127").Value = 4500
$ws.Range("L127").Value = 13500
$ws.Range("N127").Value = -23420

$ws.Range("H129").Value = 1248.75
$ws.Range("I129").Value = 565
$ws.Range("J129").Value = 1476.6666
$ws.Range("K129").Value = 1695
$ws.Range("L129").Value = 4429.9998
$ws.Range("M129").Value = 3305
$ws.Range("N129").Value = -14429.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 23548100
$ws.Range("I11").Value = 55500000
$ws.Range("J11").Value = 2246834
$ws.Range("K11").Value = 55500000
$ws.Range("L11").Value = 2246834
$ws.Range("M11").Value = -55499861
$ws.Range("N11").Value = -2247112

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1524.6666
$ws.Range("I61").Value = 1524.6666
$ws.Range("K61").Value = 1524.6666
$ws.Range("M61").Value = -1322.6666

$ws.Range("H113").Value = 1524.6666
$ws.Range("I113").Value = 1524.6666
$ws.Range("K113").Value = 1524.6666
$ws.Range("M113").Value = 645.3334

$ws.Range("H136").Value = 4032.75
$ws.Range("I136").Value = 2128.8
$ws.Range("K136").Value = 6386.400000000001
$ws.Range("M136").Value = -3836.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 71429210
$ws.Range("I100").Value = 740.1111
$ws.Range("J100").Value = 200000450
$ws.Range("K100").Value = 1480.2222
$ws.Range("L100").Value = 400000900
$ws.Range("M100").Value = -939.2221999999999
$ws.Range("N100").Value = -400001982

$ws.Range("H132").Value = 1455.2903
$ws.Range("I132").Value = 1393.6522
$ws.Range("K132").Value = 4180.9566
$ws.Range("M132").Value = -1650.9566
